# Fixed initialisation bug with Isotopes.txt
# Updates the JADE Config.xlsx default-settings workbook to point at the
# new "drake" machine paths / Slurm batch system instead of the old
# "freia" / LoadLeveler setup.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("MAIN Config.")
$wsComp = $wb.Worksheets.Item("Computational benchmarks")
$wsLib  = $wb.Worksheets.Item("Libraries")

# NOTE: the shared-string table is rebuilt on save in the order new
# strings are *written*, so the order of the .Value assignments below
# matters and mirrors the order the strings appear in the target file.

# OpenMC executable path (freia -> drake)
$wsMain.Range("B6").Value = "/home/sbradnam/Software/drake/DAG_OPENMC_021222/openmc/build/bin/openmc"

# MCNP executable path (freia -> drake)
$wsMain.Range("B2").Value = "/home/mcnp/mcnpexecs/drake/mcnp6v2_ifort2018_n1s"

# FENDL cross_sections.xml config path (CUMULUS -> NEUTRONS on drake)
$wsLib.Range("F4").Value = "/home/sbradnam/Software/NEUTRONS/openmc/data/fendl-3.1d-hdf5/cross_sections.xml"

# ENDF/B-VIII cross_sections.xml config path (CUMULUS -> NEUTRONS on drake)
$wsLib.Range("F7").Value = "/home/sbradnam/Software/NEUTRONS/openmc/data/endfb80_hdf5/cross_sections.xml"

# Batch system: llsubmit -> sbatch
$wsMain.Range("B12").Value = "sbatch"

# Batch file: LLtemplate.cmd -> Slurmtemplate.sh
$wsMain.Range("B13").Value = "Job_Script_Templates/Slurmtemplate.sh"

# MPI tasks: 4 -> 8
$wsMain.Range("B11").Value = 8

# ---------------------------------------------------------------------
# Selections
# ---------------------------------------------------------------------

# "Computational benchmarks": move the selection to F4
$wsComp.Range("F4").Select()

# "Libraries": was topLeftCell D1 + selection F3, now default top-left
# with selection D4
$wsLib.Range("A1").Select()
$wsLib.Range("D4").Select()

# "MAIN Config.": move the selection to B13 (last edited cell) and make
# sure it stays the active sheet/tab, matching the original workbook
$wsMain.Range("B13").Select()
$wsMain.Activate()
